$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-28 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-12-29 Monday", 2) | Out-Null
$d.Content.Find.Execute("316×6=1896", $true, $true, $false, $false, $false, $true, 1, $false, "933×9=8397", 2) | Out-Null
$d.Content.Find.Execute("348×3=1044", $true, $true, $false, $false, $false, $true, 1, $false, "343×9=3087", 2) | Out-Null
$d.Content.Find.Execute("147×9=1323", $true, $true, $false, $false, $false, $true, 1, $false, "937×9=8433", 2) | Out-Null
$d.Content.Find.Execute("572×2=1144", $true, $true, $false, $false, $false, $true, 1, $false, "368×4=1472", 2) | Out-Null
$d.Content.Find.Execute("231×3=693", $true, $true, $false, $false, $false, $true, 1, $false, "786×9=7074", 2) | Out-Null
$d.Content.Find.Execute("395×4=1580", $true, $true, $false, $false, $false, $true, 1, $false, "984×4=3936", 2) | Out-Null
$d.Content.Find.Execute("773×9=6957", $true, $true, $false, $false, $false, $true, 1, $false, "638×6=3828", 2) | Out-Null
$d.Content.Find.Execute("404×4=1616", $true, $true, $false, $false, $false, $true, 1, $false, "288×4=1152", 2) | Out-Null
$d.Content.Find.Execute("324×2=648", $true, $true, $false, $false, $false, $true, 1, $false, "595×2=1190", 2) | Out-Null
$d.Content.Find.Execute("339×9=3051", $true, $true, $false, $false, $false, $true, 1, $false, "238×9=2142", 2) | Out-Null
$d.Content.Find.Execute("643×4=2572", $true, $true, $false, $false, $false, $true, 1, $false, "956×4=3824", 2) | Out-Null
$d.Content.Find.Execute("599×9=5391", $true, $true, $false, $false, $false, $true, 1, $false, "299×3=897", 2) | Out-Null
$d.Content.Find.Execute("427×3=1281", $true, $true, $false, $false, $false, $true, 1, $false, "321×6=1926", 2) | Out-Null
$d.Content.Find.Execute("660×5=3300", $true, $true, $false, $false, $false, $true, 1, $false, "973×7=6811", 2) | Out-Null
$d.Content.Find.Execute("160×3=480", $true, $true, $false, $false, $false, $true, 1, $false, "917×8=7336", 2) | Out-Null
$d.Content.Find.Execute("621×6=3726", $true, $true, $false, $false, $false, $true, 1, $false, "941×6=5646", 2) | Out-Null
$d.Content.Find.Execute("730×6=4380", $true, $true, $false, $false, $false, $true, 1, $false, "337×7=2359", 2) | Out-Null
$d.Content.Find.Execute("211×2=422", $true, $true, $false, $false, $false, $true, 1, $false, "143×8=1144", 2) | Out-Null
$d.Content.Find.Execute("308×9=2772", $true, $true, $false, $false, $false, $true, 1, $false, "451×8=3608", 2) | Out-Null
$d.Content.Find.Execute("807×9=7263", $true, $true, $false, $false, $false, $true, 1, $false, "435×9=3915", 2) | Out-Null
$d.Content.Find.Execute("352×9=3168", $true, $true, $false, $false, $false, $true, 1, $false, "914×6=5484", 2) | Out-Null
$d.Content.Find.Execute("919×6=5514", $true, $true, $false, $false, $false, $true, 1, $false, "359×3=1077", 2) | Out-Null
$d.Content.Find.Execute("550×7=3850", $true, $true, $false, $false, $false, $true, 1, $false, "980×9=8820", 2) | Out-Null
$d.Content.Find.Execute("533×5=2665", $true, $true, $false, $false, $false, $true, 1, $false, "578×6=3468", 2) | Out-Null
$d.Content.Find.Execute("823×6=4938", $true, $true, $false, $false, $false, $true, 1, $false, "523×5=2615", 2) | Out-Null
